$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values stay as text (preserve exact formatting, trailing zeros, % signs)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "307.11"
$ws.Range("E2").Value = "1.13%"
$ws.Range("D3").Value = "36.99"
$ws.Range("E3").Value = "3.73%"
$ws.Range("D4").Value = "5.043"
$ws.Range("E4").Value = "-0.63%"
$ws.Range("D5").Value = "0.07888"
$ws.Range("E5").Value = "1.09%"
$ws.Range("D6").Value = "2.192"
$ws.Range("E6").Value = "-2.99%"
$ws.Range("D7").Value = "8.041"
$ws.Range("E7").Value = "-0.69%"
$ws.Range("D8").Value = "0.9273"
$ws.Range("E8").Value = "-0.29%"
$ws.Range("D9").Value = "0.09945"
$ws.Range("E9").Value = "2.21%"
$ws.Range("D10").Value = "0.1891"
$ws.Range("E10").Value = "3.76%"
$ws.Range("D11").Value = "0.08723"
$ws.Range("E11").Value = "-0.21%"
$ws.Range("D12").Value = "0.03600"
$ws.Range("E12").Value = "5.09%"
$ws.Range("D13").Value = "0.09965"
$ws.Range("E13").Value = "0.41%"
$ws.Range("D14").Value = "0.001497"
$ws.Range("E14").Value = "1.69%"
$ws.Range("D15").Value = "0.005662"
$ws.Range("E15").Value = "0.49%"
$ws.Range("D16").Value = "3.463"
$ws.Range("E16").Value = "-0.80%"
$ws.Range("D17").Value = "4.058"
$ws.Range("E17").Value = "1.19%"
$ws.Range("D18").Value = "2.327"
$ws.Range("E18").Value = "9.42%"
$ws.Range("D19").Value = "0.3435"
$ws.Range("E19").Value = "0.13%"
$ws.Range("E20").Value = "0.38%"
$ws.Range("D21").Value = "4.921"
$ws.Range("E21").Value = "7.89%"
$ws.Range("D22").Value = "0.2200"
$ws.Range("E22").Value = "-1.71%"
$ws.Range("D23").Value = "0.04636"
$ws.Range("E23").Value = "-0.84%"
$ws.Range("D24").Value = "0.005230"
$ws.Range("E24").Value = "16.50%"
$ws.Range("D25").Value = "0.001249"
$ws.Range("E25").Value = "0.76%"
$ws.Range("D26").Value = "0.0001400"
$ws.Range("E26").Value = "7.51%"
$ws.Range("D27").Value = "0.0002716"
$ws.Range("E27").Value = "0.52%"
$ws.Range("D39").Value = "0.01837"
$ws.Range("E39").Value = "4.71%"
$ws.Range("D40").Value = "0.04778"
$ws.Range("E40").Value = "1.67%"
$ws.Range("D41").Value = "0.007966"
$ws.Range("E41").Value = "1.44%"
$ws.Range("D42").Value = "0.1417"
$ws.Range("E42").Value = "-0.05%"
$ws.Range("D43").Value = "0.007586"
$ws.Range("E43").Value = "-12.28%"
$ws.Range("D44").Value = "0.002190"
$ws.Range("E44").Value = "-4.53%"
$ws.Range("D45").Value = "0.01010"
$ws.Range("E45").Value = "9.57%"
$ws.Range("D46").Value = "0.00006262"
$ws.Range("E46").Value = "2.07%"
$ws.Range("E47").Value = "-0.15%"
$ws.Range("D48").Value = "0.0005799"
$ws.Range("E48").Value = "-0.02%"
$ws.Range("D49").Value = "32.21"
$ws.Range("E49").Value = "457.29%"
$ws.Range("D50").Value = "0.002688"
$ws.Range("E50").Value = "-0.18%"
$ws.Range("D51").Value = "0.00002100"
$ws.Range("E51").Value = "-0.15%"

# Coin name / link updates (plain text, no special formatting needed)
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
